$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '27.999.03'
$ws.Range('E2').Value = '  +1.62%  '
$ws.Range('D3').Value = '1.755.30'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.82'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3836'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3423'
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.77'
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.123'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07241'
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.57'
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.177'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.143'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '1.750.25'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001063'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06612'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '79.39'
$ws.Range('E19').Value = '  -3.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9993'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.78'
$ws.Range('E21').Value = '  -3.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.201'
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('D23').Value = '27.998.66'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.69'
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.380'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.95'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('E27').Value = '  -3.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.305'
$ws.Range('E28').Value = '  -4.47%  '
$ws.Range('D29').Value = '1.950.38'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.262'
$ws.Range('E30').Value = '  -11.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '131.15'
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.032'
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.861'
$ws.Range('E33').Value = '  -3.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08822'
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.22'
$ws.Range('E35').Value = '  -3.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.550'
$ws.Range('E36').Value = '  +3.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6582'
$ws.Range('E37').Value = '  -2.46%  '
$ws.Range('E38').Value = '  -4.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.149'
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06161'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2108'
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.217'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.967'
$ws.Range('E43').Value = '  -3.15%  '
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.73'
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.835'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6054'
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.35'
$ws.Range('E48').Value = '  -2.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.007'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('E50').Value = '  +2.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.112'
$ws.Range('E51').Value = '  +4.93%  '
